# Update "想去人数" (want-to-go count) figures on the 展览 (exhibition)
# and 全部类型 (all types) sheets to the latest scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 369
$ws1.Range("F5").Value = 417
$ws1.Range("F6").Value = 264
$ws1.Range("F7").Value = 2412
$ws1.Range("F8").Value = 415
$ws1.Range("F9").Value = 6284
$ws1.Range("F11").Value = 405

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 369
$ws4.Range("F5").Value = 417
$ws4.Range("F6").Value = 264
$ws4.Range("F9").Value = 2412
$ws4.Range("F10").Value = 415
$ws4.Range("F11").Value = 6284
$ws4.Range("F13").Value = 405
